$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "5 - Make Reservation"
$ws.Range("A6").Value = "6 - Change Reservation"
$ws.Range("A7").Value = "7 - Remove/Cancel  Reservation"

$ws.Range("A6:A7").HorizontalAlignment = -4131

$ws.Range("B12").Select()
